$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.334811329841614
$ws.Range("B1").Value = 2.940743207931519
$ws.Range("C1").Value = 1.822762727737427
$ws.Range("D1").Value = 1.261430263519287
$ws.Range("E1").Value = 1.045578598976135
